$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every Price (D) / Volume(1h) (E) cell we touch, force Text format
# first so Excel does not auto-convert numeric-looking strings (e.g. "1.003",
# "27.657.73", "47.70") into numbers - the source stores these as literal text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.657.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.70"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08043"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.86"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.879.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.949"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.140"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.663.50"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.092.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.550"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.23"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9689"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09500"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.445"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.587"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.324"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02259"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06090"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.122"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5989"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1893"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5678"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.23"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.388"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06848"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.33"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.71%  "
